$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.786.83'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.14%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.118.77'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +10.36%  '

$ws.Range('E4').Value = '  -0.15%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '334.62'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.65%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9996'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.17%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5247'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.63%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4413'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +8.32%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09104'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +9.00%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '47.15'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +11.33%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.186'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +6.70%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '25.32'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +5.24%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.119.62'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +10.71%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.791'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.60%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.858'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +8.33%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '98.18'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +6.03%  '

$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001140'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.03%  '

$ws.Range('B18').Value = 'BinanceUSD'
$ws.Range('C18').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.000'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.27%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06652'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.05%  '

$ws.Range('E20').Value = '  +3.71%  '

$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.412'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +7.78%  '

$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9997'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.19%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '30.902.90'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.50%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.10'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +6.40%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.364.41'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +10.80%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.255'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.78%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.98'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.94%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.563'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +12.97%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '163.51'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.59%  '

$ws.Range('E30').Value = '  +3.64%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.186'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.65%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1073'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.62%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.274'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.27%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.964'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.66%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.539'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +28.41%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02605'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.06%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.602'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.09%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '9.620'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +11.97%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06776'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.20%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '12.79'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +11.77%  '

$ws.Range('E41').Value = '  +5.70%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6849'
$ws.Range('D42').Style = 'Normal'

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.257'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.92%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.25'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.09%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6436'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.15%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9992'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.14%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.269'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.31%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.679'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.47%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.288'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +6.30%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '83.28'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.44%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07090'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.85%  '
